$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4134
$ws.Range("I18").Value = 3333.3333
$ws.Range("J18").Value = 4934.6665
$ws.Range("K18").Value = 3333.3333
$ws.Range("L18").Value = 4934.6665
$ws.Range("M18").Value = -3049.3333
$ws.Range("N18").Value = -5502.6665
$ws.Range("H21").Value = 7509.5
$ws.Range("J21").Value = 7509.5
$ws.Range("L21").Value = 7509.5
$ws.Range("N21").Value = -8445.5
$ws.Range("H23").Value = 7509.5
$ws.Range("J23").Value = 7509.5
$ws.Range("L23").Value = 7509.5
$ws.Range("N23").Value = -7977.5
$ws.Range("H40").Value = 999
$ws.Range("I40").Value = 999
$ws.Range("K40").Value = 999
$ws.Range("M40").Value = -824
$ws.Range("H51").Value = 6666
$ws.Range("I51").Value = 6666
$ws.Range("K51").Value = 6666
$ws.Range("M51").Value = -6182
$ws.Range("H58").Value = 1493.5834
$ws.Range("J58").Value = 3348.75
$ws.Range("L58").Value = 10046.25
$ws.Range("N58").Value = -10346.25
$ws.Range("H82").Value = 1183.5714
$ws.Range("I82").Value = 1183.5714
$ws.Range("K82").Value = 3550.7142
$ws.Range("M82").Value = -3144.7142
$ws.Range("H85").Value = 1183.5714
$ws.Range("I85").Value = 1183.5714
$ws.Range("K85").Value = 3550.7142
$ws.Range("M85").Value = -2146.7142
$ws.Range("H86").Value = 5098.8
$ws.Range("I86").Value = 5164.8335
$ws.Range("J86").Value = 4999.75
$ws.Range("K86").Value = 5164.8335
$ws.Range("L86").Value = 4999.75
$ws.Range("M86").Value = -4041.8335
$ws.Range("N86").Value = -7245.75
$ws.Range("H87").Value = 38749.25
$ws.Range("J87").Value = 38749.25
$ws.Range("L87").Value = 38749.25
$ws.Range("N87").Value = -41245.25
$ws.Range("H89").Value = 5098.8
$ws.Range("I89").Value = 5164.8335
$ws.Range("J89").Value = 4999.75
$ws.Range("K89").Value = 25824.1675
$ws.Range("L89").Value = 24998.75
$ws.Range("M89").Value = -20208.1675
$ws.Range("N89").Value = -36230.75
$ws.Range("H90").Value = 38749.25
$ws.Range("J90").Value = 38749.25
$ws.Range("L90").Value = 116247.75
$ws.Range("N90").Value = -128727.75
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2320.2666
$ws.Range("I61").Value = 1950.2858
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 1950.2858
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -1738.2858
$ws.Range("N61").Value = -7924
$ws.Range("H110").Value = 2990.4546
$ws.Range("I110").Value = 579
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 579
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = 1466
$ws.Range("N110").Value = -9090
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1633.3334
$ws.Range("K122").Value = 4900.0002
$ws.Range("M122").Value = -2450.0002
$ws.Range("H136").Value = 2320.2666
$ws.Range("I136").Value = 1950.2858
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 5850.857400000001
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -3300.857400000001
$ws.Range("N136").Value = -27600

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 659.9091
$ws.Range("J107").Value = 420
$ws.Range("L107").Value = 420
$ws.Range("N107").Value = -4260
$ws.Range("H134").Value = 6625.357
$ws.Range("I134").Value = 6625.357
$ws.Range("K134").Value = 19876.071
$ws.Range("M134").Value = -17341.071

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 16399.6
$ws.Range("J14").Value = 16399.6
$ws.Range("L14").Value = 16399.6
$ws.Range("N14").Value = -16739.6
$ws.Range("H41").Value = 18600
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20856
$ws.Range("H50").Value = 29998
$ws.Range("J50").Value = 29998
$ws.Range("L50").Value = 29998
$ws.Range("N50").Value = -31248
$ws.Range("H58").Value = 2027.2727
$ws.Range("I58").Value = 1737.5
$ws.Range("K58").Value = 1737.5
$ws.Range("M58").Value = -1534.5
$ws.Range("H60").Value = 16341.857
$ws.Range("H68").Value = 38588.43
$ws.Range("J68").Value = 39997.69
$ws.Range("L68").Value = 39997.69
$ws.Range("N68").Value = -41495.69
$ws.Range("H71").Value = 38588.43
$ws.Range("J71").Value = 39997.69
$ws.Range("L71").Value = 119993.07
$ws.Range("N71").Value = -127481.07
$ws.Range("H74").Value = 38081
$ws.Range("J74").Value = 38225.273
$ws.Range("L74").Value = 38225.273
$ws.Range("N74").Value = -39973.273
$ws.Range("H77").Value = 38081
$ws.Range("J77").Value = 38225.273
$ws.Range("L77").Value = 114675.819
$ws.Range("N77").Value = -123411.819
$ws.Range("H136").Value = 2027.2727
$ws.Range("I136").Value = 1737.5
$ws.Range("K136").Value = 5212.5
$ws.Range("M136").Value = -2662.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1267.7142
$ws.Range("J4").Value = 1333.1765
$ws.Range("L4").Value = 3999.5295
$ws.Range("N4").Value = -4223.529500000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1526
$ws.Range("I102").Value = 1526
$ws.Range("K102").Value = 1526
$ws.Range("M102").Value = 96
$ws.Range("H122").Value = 11368109
$ws.Range("I122").Value = 13893122
$ws.Range("K122").Value = 41679366
$ws.Range("M122").Value = -41676916

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2973.2222
$ws.Range("I46").Value = 2685.7144
$ws.Range("K46").Value = 2685.7144
$ws.Range("M46").Value = -2497.7144
$ws.Range("H55").Value = 180.5
$ws.Range("I55").Value = 136.14285
$ws.Range("J55").Value = 284
$ws.Range("K55").Value = 136.14285
$ws.Range("L55").Value = 284
$ws.Range("M55").Value = 36.85714999999999
$ws.Range("N55").Value = -630
$ws.Range("H61").Value = 7921.8
$ws.Range("I61").Value = 8024.222
$ws.Range("K61").Value = 8024.222
$ws.Range("M61").Value = -7822.222
$ws.Range("H100").Value = 2999
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458
$ws.Range("H113").Value = 7921.8
$ws.Range("I113").Value = 8024.222
$ws.Range("K113").Value = 8024.222
$ws.Range("M113").Value = -5854.222

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 597.1429000000001
$ws.Range("I107").Value = 597.1429000000001
$ws.Range("K107").Value = 1791.4287
$ws.Range("M107").Value = 128.5712999999998
$ws.Range("H113").Value = 955.6667
$ws.Range("I113").Value = 1016.8333
$ws.Range("K113").Value = 3050.4999
$ws.Range("M113").Value = -880.4998999999998
$ws.Range("H122").Value = 1776.5714
$ws.Range("I122").Value = 1592.909
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 4778.727000000001
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -2328.727000000001
$ws.Range("N122").Value = -12250

